# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Elimina los trabajadores "JOSE ANDRES CASTILLA ROMERO" (1100393071) y
#   "LUIS ANGEL HURTADO RAMIREZ" (1193599222), cuyos valores de mora no
#   coincidian con el resto del periodo 2507.
# - Agrega un segundo bloque (periodo 2508) con los 6 trabajadores que
#   quedan, replicando su valor de mora / salario basico.
# - Actualiza los totales de "VALOR MORA" y "Cant. Trabajadores".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Quitar los dos trabajadores que salen de la base -----------------
# Fila 17 = JOSE ANDRES CASTILLA ROMERO (1100393071)
$ws.Rows.Item(17).Delete()
# Tras el borrado anterior, LUIS ANGEL HURTADO RAMIREZ (1193599222) queda
# en la fila 21 (antes fila 22)
$ws.Rows.Item(21).Delete()

# En este punto quedan 6 trabajadores del periodo 2507 en las filas 16-21,
# con DAVINZON CABALLERO ARELLANO en la fila 21 (ultima fila, con el borde
# inferior grueso propio del cierre de la tabla).

# --- 2. Insertar 6 filas nuevas para el periodo 2508 ----------------------
$ws.Range("22:27").Insert()

# Conservar el formato de "ultima fila" (borde inferior grueso) copiandolo
# a la nueva ultima fila (27) antes de normalizar la fila 21.
$ws.Range("B21:J21").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)

# Ahora la fila 21 vuelve a ser una fila "normal" de la tabla, igual que
# las filas nuevas 22-26: copiar el formato de la fila 20.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J26").PasteSpecial(-4122)

# --- 3. Rellenar los datos del periodo 2508 --------------------------------
$employees = @(
    @("CC", "73201709",   "ULEY RODRIGUEZ ORTEGA"),
    @("CC", "9286755",    "WILLIAM RAFAEL NAVARRO PEREZ"),
    @("CC", "12536656",   "JAIME ENRIQUE PACHECO SOTO"),
    @("CC", "19611421",   "JOAQUIN ANTONIO TAPIAS BORNACHERA"),
    @("CC", "1047512837", "DARIN ESCORCIA CAMARGO"),
    @("CC", "9296960",    "DAVINZON CABALLERO ARELLANO")
)

$r = 22
foreach ($emp in $employees) {
    $ws.Cells.Item($r, 2).Value = $emp[0]
    $ws.Cells.Item($r, 3).Value = $emp[1]
    $ws.Cells.Item($r, 4).Value = $emp[2]
    $ws.Cells.Item($r, 5).Value = "2508"
    $ws.Cells.Item($r, 6).Value = 56940
    $ws.Cells.Item($r, 7).Value = 1423500
    $r = $r + 1
}

# --- 4. Actualizar totales --------------------------------------------------
# Cant. Trabajadores: 8 -> 6
$ws.Range("C13").Value = 6
# VALOR MORA total: 399334 -> 683280 (6 trabajadores x 2 periodos x 56940)
$ws.Range("E11").Value = 683280
